# Update column G ("K") values per the regenerated save_data.
# The author's commit regenerated the save_data workbook to use K (strikeouts)
# instead of Strike# in column G for rows 2-75 (data rows).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newK = @{
    2  = 0
    3  = 1
    4  = 2
    5  = 0
    6  = 1
    7  = 1
    8  = 0
    9  = 1
    10 = 3
    11 = 0
    12 = 0
    13 = 1
    14 = 1
    15 = 1
    16 = 0
    17 = 1
    18 = 1
    19 = 0
    20 = 1
    21 = 0
    22 = 1
    23 = 0
    24 = 0
    25 = 0
    26 = 0
    27 = 0
    28 = 1
    29 = 2
    30 = 0
    31 = 2
    32 = 2
    33 = 1
    34 = 1
    35 = 0
    37 = 1
    38 = 0
    39 = 1
    40 = 2
    41 = 3
    42 = 2
    43 = 0
    44 = 2
    45 = 1
    46 = 1
    47 = 0
    48 = 0
    49 = 0
    50 = 2
    51 = 0
    52 = 1
    53 = 0
    54 = 2
    55 = 1
    56 = 0
    57 = 0
    58 = 0
    59 = 1
    60 = 0
    61 = 2
    62 = 1
    63 = 1
    64 = 0
    65 = 0
    66 = 1
    67 = 0
    68 = 2
    69 = 2
    70 = 2
    71 = 2
    72 = 2
    75 = 0
}

foreach ($row in $newK.Keys) {
    $ws.Range("G$row").Value = $newK[$row]
}
